$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style from E1 (existing header) to F1:H1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean values for rows 2-8
$values = @(
    @($true, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($true, $true, $true),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $values[$i][0]
    $ws.Range("G$row").Value = $values[$i][1]
    $ws.Range("H$row").Value = $values[$i][2]
}
